$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list: updated Price (col D) and 1h Volume change (col E)
# for every coin row (2-51), matching the latest scrape from the GitHub Action.
#
# Price strings that look like plain numbers (e.g. "328.40") are written with a
# leading apostrophe so Excel stores them as literal text (preserving trailing
# zeros / exact formatting) instead of silently parsing them into a float and
# losing the original "328.40" -> 328.4 formatting. Prices that already contain
# two dots (e.g. "27.371.84") are never valid numbers, so they are assigned directly.

$ws.Range('D2').Value = '27.371.84'
$ws.Range('E2').Value = '  -2.90%  '
$ws.Range('D3').Value = '1.857.89'
$ws.Range('E3').Value = '  -3.00%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''328.40'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '''0.4576'
$ws.Range('E7').Value = '  -2.10%  '
$ws.Range('D8').Value = '''0.3929'
$ws.Range('E8').Value = '  -1.92%  '
$ws.Range('D9').Value = '''47.13'
$ws.Range('E9').Value = '  -11.29%  '
$ws.Range('D10').Value = '''0.07939'
$ws.Range('E10').Value = '  -5.46%  '
$ws.Range('D11').Value = '''1.013'
$ws.Range('E11').Value = '  -2.93%  '
$ws.Range('E12').Value = '  -2.77%  '
$ws.Range('D13').Value = '1.871.64'
$ws.Range('E13').Value = '  -2.51%  '
$ws.Range('D14').Value = '''5.923'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').Value = '''7.151'
$ws.Range('E15').Value = '  -3.51%  '
$ws.Range('D16').Value = '''1.003'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').Value = '''86.45'
$ws.Range('E17').Value = '  -3.49%  '
$ws.Range('D18').Value = '''0.06609'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('D19').Value = '''0.00001031'
$ws.Range('E19').Value = '  -2.81%  '
$ws.Range('D20').Value = '''17.24'
$ws.Range('E20').Value = '  -3.95%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '''5.488'
$ws.Range('E22').Value = '  -4.00%  '
$ws.Range('D23').Value = '27.368.91'
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('D24').Value = '''10.93'
$ws.Range('E24').Value = '  -3.34%  '
$ws.Range('D25').Value = '''2.306'
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').Value = '2.073.19'
$ws.Range('E26').Value = '  -3.35%  '
$ws.Range('D27').Value = '''153.66'
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').Value = '''20.08'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('D29').Value = '''2.070'
$ws.Range('E29').Value = '  -2.53%  '
$ws.Range('D30').Value = '''5.470'
$ws.Range('E30').Value = '  -4.19%  '
$ws.Range('D31').Value = '''121.70'
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').Value = '''0.9534'
$ws.Range('E32').Value = '  -1.97%  '
$ws.Range('D33').Value = '''0.09402'
$ws.Range('E33').Value = '  -1.95%  '
$ws.Range('D34').Value = '''1.458'
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('D35').Value = '''3.589'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('D36').Value = '''5.275'
$ws.Range('E36').Value = '  -4.74%  '
$ws.Range('D37').Value = '''0.06032'
$ws.Range('E37').Value = '  -2.09%  '
$ws.Range('D38').Value = '''0.02230'
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('D39').Value = '''1.217'
$ws.Range('E39').Value = '  -2.18%  '
$ws.Range('D40').Value = '''8.052'
$ws.Range('E40').Value = '  -8.39%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').Value = '''0.5923'
$ws.Range('E42').Value = '  -3.23%  '
$ws.Range('D43').Value = '''0.1886'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('E44').Value = '  -7.56%  '
$ws.Range('D45').Value = '''1.281'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').Value = '''0.5620'
$ws.Range('E46').Value = '  -3.65%  '
$ws.Range('D47').Value = '''12.11'
$ws.Range('E47').Value = '  -4.58%  '
$ws.Range('D48').Value = '''3.397'
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('D49').Value = '''1.919'
$ws.Range('E49').Value = '  -5.11%  '
$ws.Range('D50').Value = '''0.06749'
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('D51').Value = '''108.30'
$ws.Range('E51').Value = '  -1.45%  '
